# Add team record columns (Wins / Losses / Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers in AD1:AF1, matching the existing header style (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-44: constant team record values.
$ws.Range("AD2:AD44").Value = 90
$ws.Range("AE2:AE44").Value = 72
$ws.Range("AF2:AF44").Value = 0
